# ENH: Parameterize assumption file prefix
#
# On the ParamList sheet, rename the existing "model_point_file_prefix"
# parameter row to "mp_file_prefix" and insert a new row right below it
# for a new "asmp_file_prefix" parameter. Mirror the same insertion on
# the ConstParams sheet (which holds the actual constant values), and
# update the active sheet / selections to match.

$wb = $excel.ActiveWorkbook

# --- ParamList sheet ---------------------------------------------------
$paramList = $wb.Worksheets.Item("ParamList")
$constParams = $wb.Worksheets.Item("ConstParams")

# Row 4 currently holds "model_point_file_prefix" / CONST / "Stem part
# of the model point file name". Insert a fresh row underneath it for
# the new "asmp_file_prefix" parameter, fill in its contents (mirroring
# the new constant row on ConstParams at the same time), then finally
# rename the original row's name to "mp_file_prefix".
$paramList.Rows.Item(5).Insert()
$constParams.Rows.Item(4).Insert()

$paramList.Cells.Item(5, 1).Value = "asmp_file_prefix"
$paramList.Cells.Item(5, 2).Value = "CONST"

$constParams.Cells.Item(4, 1).Value = "asmp_file_prefix"
$constParams.Cells.Item(4, 2).Value = "assumptions"

$paramList.Cells.Item(5, 3).Value = "Stem part of the assumption file name"

$paramList.Cells.Item(4, 1).Value = "mp_file_prefix"
$constParams.Cells.Item(3, 1).Value = "mp_file_prefix"

# The autofilter only ever covered the parameter rows (A1:C15); it now
# needs to cover one extra row (A1:C16) after the insertion above, while
# leaving the later, unfiltered rows (17-19) alone. The host always grows
# a freshly-(re)applied AutoFilter down to the full contiguous block of
# data below the header, so temporarily blank out those trailing rows,
# apply the filter over exactly the rows it should cover, then restore
# the stashed values — the filter's ref stays put once it's been set.
$row17 = @($paramList.Cells.Item(17, 1).Value2, $paramList.Cells.Item(17, 2).Value2, $paramList.Cells.Item(17, 3).Value2)
$row18 = @($paramList.Cells.Item(18, 1).Value2, $paramList.Cells.Item(18, 2).Value2, $paramList.Cells.Item(18, 3).Value2)
$row19 = @($paramList.Cells.Item(19, 1).Value2, $paramList.Cells.Item(19, 2).Value2, $paramList.Cells.Item(19, 3).Value2)

$paramList.AutoFilterMode = $false
$paramList.Rows.Item(17).ClearContents()
$paramList.Rows.Item(18).ClearContents()
$paramList.Rows.Item(19).ClearContents()

$paramList.Range("A1:C16").AutoFilter()

$paramList.Cells.Item(17, 1).Value = $row17[0]
$paramList.Cells.Item(17, 2).Value = $row17[1]
$paramList.Cells.Item(17, 3).Value = $row17[2]
$paramList.Cells.Item(18, 1).Value = $row18[0]
$paramList.Cells.Item(18, 2).Value = $row18[1]
$paramList.Cells.Item(18, 3).Value = $row18[2]
$paramList.Cells.Item(19, 1).Value = $row19[0]
$paramList.Cells.Item(19, 2).Value = $row19[1]
$paramList.Cells.Item(19, 3).Value = $row19[2]

# Keep the workbook-level _FilterDatabase defined name in sync with the
# resized autofilter range.
$wb.Names.Item(1).RefersTo = "=ParamList!`$A`$1:`$C`$16"

$paramList.Range("C20").Select()
$constParams.Range("B14").Select()

# --- Active sheet / selections on other sheets --------------------------
$runParams = $wb.Worksheets.Item("RunParams")
$runParams.Range("E28").Select()

$gmxb = $wb.Worksheets.Item("GMXB")
$gmxb.Range("H32").Select()

$constParams.Activate()
